# Updated list to reflect latest apps versions from web
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update app version numbers (column C) ---
$ws.Range("C2").Value  = "0.31.0"     # visual_c
$ws.Range("C5").Value  = "4.19.0"     # calibre
$ws.Range("C7").Value  = "8.6.1"      # crystaldiskinfo
$ws.Range("C9").Value  = "99.4.501"   # dropbox
$ws.Range("C13").Value = "12.00"      # exiftool
$ws.Range("C15").Value = "77.0.1"     # firefox
$ws.Range("C16").Value = "2.27.0"     # git
$ws.Range("C20").Value = "12.10.7"    # itunes
$ws.Range("C21").Value = "15.5.0"     # klite_codec
$ws.Range("C22").Value = "47.0.0"     # mkvtoolnix
$ws.Range("C26").Value = "3.11"       # rufus
$ws.Range("C28").Value = "1.46"       # visual_studio_code
$ws.Range("C29").Value = "5.17.6"     # winscp

# --- Split the default-width column block so column 14 becomes its own
#     explicit entry (matching columns 7-13), leaving 15-16384 grouped ---
$ws.Columns.Item(14).Hidden = $False

# --- Move the sheet selection from A1 to D31 ---
$ws.Range("D31").Select()
